# Update COVID-19 country statistics and the "last updated" timestamp
# (values refreshed to the 17:20 snapshot of 3 April 2020).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1: "Datos actualizados a 3 de Abril de 2020 a las 16:50" -> "...17:20"
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 17:20"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 257379
$ws.Range("C4").Value = 12502
$ws.Range("D4").Value = 11941
$ws.Range("E4").Value = 238880
$ws.Range("F4").Value = 5781
$ws.Range("G4").Value = 488
$ws.Range("H4").Value = 6558

# Row 7: Alemania
$ws.Range("B7").Value = 89126
$ws.Range("C7").Value = 4332
$ws.Range("E7").Value = 63353
$ws.Range("G7").Value = 91
$ws.Range("H7").Value = 1198

# Row 16: Austria
$ws.Range("B16").Value = 11444
$ws.Range("C16").Value = 315
$ws.Range("E16").Value = 9254

# Row 42: Finlandia
$ws.Range("F42").Value = 72

# Row 43: Grecia
$ws.Range("B43").Value = 1613
$ws.Range("C43").Value = 69
$ws.Range("E43").Value = 1493
$ws.Range("F43").Value = 92
$ws.Range("G43").Value = 6
$ws.Range("H43").Value = 59

# Row 45: Republica Dominicana
$ws.Range("B45").Value = 1505
$ws.Range("C45").Value = 43
$ws.Range("D45").Value = 95
$ws.Range("E45").Value = 1403
$ws.Range("F45").Value = 7
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 7

# Row 46: Serbia
$ws.Range("B46").Value = 1488
$ws.Range("C46").Value = 108
$ws.Range("D46").Value = 16
$ws.Range("E46").Value = 1404
$ws.Range("F46").Value = 147
$ws.Range("H46").Value = 68

# Row 47: Panama
$ws.Range("B47").Value = 1476
$ws.Range("C47").Value = 305
$ws.Range("D47").Value = 42
$ws.Range("E47").Value = 1395
$ws.Range("F47").Value = 81
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 39

# Row 48: Sudafrica
$ws.Range("B48").Value = 1475
$ws.Range("D48").Value = 9
$ws.Range("E48").Value = 1429
$ws.Range("F48").Value = 50
$ws.Range("H48").Value = 37

# Row 64: Irak
$ws.Range("B64").Value = 820
$ws.Range("C64").Value = 48
$ws.Range("D64").Value = 226
$ws.Range("E64").Value = 540

# Row 84: Uruguay
$ws.Range("B84").Value = 396
$ws.Range("C84").Value = 40
$ws.Range("D84").Value = 28
$ws.Range("E84").Value = 357
$ws.Range("F84").Value = 11
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 11

# Row 85: Republica de Chipre
$ws.Range("B85").Value = 369
$ws.Range("C85").Value = 19
$ws.Range("D85").Value = 68
$ws.Range("E85").Value = 297
$ws.Range("F85").Value = 13
$ws.Range("H85").Value = 4

# Row 124: Paraguay
$ws.Range("B124").Value = 95
$ws.Range("C124").Value = 7
$ws.Range("D124").Value = 46
$ws.Range("E124").Value = 49
$ws.Range("F124").Value = 0
$ws.Range("H124").Value = 0

# Row 125: Gibraltar
$ws.Range("B125").Value = 92
$ws.Range("C125").Value = 15
$ws.Range("D125").Value = 4
$ws.Range("E125").Value = 85
$ws.Range("F125").Value = 4
$ws.Range("H125").Value = 3

# Row 142: Puerto Rico
$ws.Range("D142").Value = 2
$ws.Range("H142").Value = 1

# Row 143: Zambia
$ws.Range("D143").Value = 1
$ws.Range("H143").Value = 2

# Row 160: Nueva Caledonia
$ws.Range("C160").Value = 2

# Row 161: Haiti
$ws.Range("C161").Value = 0

# Row 178: Suazilandia
$ws.Range("B178").Value = 10
$ws.Range("C178").Value = 2
$ws.Range("D178").Value = 2
$ws.Range("E178").Value = 6
$ws.Range("H178").Value = 2

# Row 182: Zimbabue
$ws.Range("E182").Value = 9
$ws.Range("H182").Value = 0

# Row 183: Montserrat
$ws.Range("E183").Value = 8
$ws.Range("H183").Value = 1

# Row 184: Republica de Africa Central
$ws.Range("B184").Value = 9
$ws.Range("C184").Value = 0
$ws.Range("E184").Value = 7
$ws.Range("H184").Value = 2

# Row 186: Angola
$ws.Range("C186").Value = 5
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 8
$ws.Range("H186").Value = 0

# Row 187: Sudan
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = 5
